$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "gameWidget" box spanning B4:B7 (same plain centered look as the
#    existing B2:B3 "selectactorwidget" box) and the new "stopGame" colored
#    slot label that lives in that row's C column.
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "gameWidget"

# ---------------------------------------------------------------------------
# 2) Column C gets a brand-new "slot name" style: same teal font used in
#    column D (fontId 1 / FF00677C Arial Unicode MS) but without the
#    center/middle alignment that column D uses. Build it once on C2, then
#    fan the formatting out to the rest of the column via copy/paste so we
#    do not re-mint a fresh font object for every single cell.
# ---------------------------------------------------------------------------
$slotFont = $ws.Range("C2")
$slotFont.ClearFormats()
$slotFont.Font.Color = 8152832
$slotFont.Font.Size = 10
$slotFont.Font.Name = "Arial Unicode MS"

$slotFont.Copy()
$ws.Range("C3:C8").PasteSpecial(-4122)

$ws.Range("C2").Value = "sendCurrentActor"
$ws.Range("C3").Value = "toGameWidget"
$ws.Range("C4").Value = "toStatsWidget"
$ws.Range("C5").Value = "toSelectActorWidget"
$ws.Range("C6").Value = "startGame"
$ws.Range("C7").Value = "stopGame"
$ws.Range("C8").Value = "toGameWidget"

# ---------------------------------------------------------------------------
# 3) New row 8: "statwidget" box in column B (plain, unstyled cell).
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "statwidget"

# ---------------------------------------------------------------------------
# 4) Merge the new gameWidget box and move the selection onto it, matching
#    the author's last interaction before saving.
# ---------------------------------------------------------------------------
$ws.Range("B4:B7").Merge()
$ws.Range("B4:B7").Select()
